# Append the new metrics row (row 30) that was logged by the newly-added
# Streamlit app, and drop the now-empty "Test R2" (F29) placeholder cell
# from the previous row the same way the upstream export script does.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The former last row (29) no longer carries an empty inline-string cell
# in column F ("Test R2") - remove it so the cell is gone entirely.
$ws.Range("F29").ClearContents()

# New row 30 - mirrors the layout of the existing rows:
# A=Date, B=Train R2, C=Train MAE, D=Train MSE, E=Train RMSE,
# F=Test R2 (left blank for this run), G=Test MAE, H=Test MSE, I=Test RMSE
$ws.Range("A30").Value = "2024-12-09 09:16:46"
$ws.Range("B30").Value = 0.9961636085978972
$ws.Range("C30").Value = 0.00762124003657671
$ws.Range("D30").Value = 0.0001757898104453401
$ws.Range("E30").Value = 0.01325857497792806
$ws.Range("G30").Value = 0.0130750132779512
$ws.Range("H30").Value = 0.0002447681723678711
$ws.Range("I30").Value = 0.01564506862777761
